$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update frequency values in column C
$ws.Range("C2").Value = 3023
$ws.Range("C3").Value = 2869
$ws.Range("C4").Value = 2136
$ws.Range("C5").Value = 1313
$ws.Range("C6").Value = 1190
$ws.Range("C7").Value = 679
$ws.Range("C8").Value = 599
$ws.Range("C9").Value = 438
$ws.Range("C10").Value = 412
$ws.Range("C11").Value = 412

# Swap category labels in B10/B11
$ws.Range("B10").Value = "Seasonal & Holidays"
$ws.Range("B11").Value = "Kitchen & Dining"
